# Added user groups changing functionality
# Adds a new "userGroupAnalysis" worksheet (after "Roles") listing the
# tourist-group taxonomy used to analyse/assign user groups.

$wb = $excel.ActiveWorkbook

# --- create the new sheet as the LAST tab -----------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$ws.Name = "userGroupAnalysis"

# --- header row (A1/C1 reuse the existing "name"/"description" strings,
#     B1 "nameForUser" is written last, matching original authoring order) -
$ws.Cells.Item(1, 1).Value = "name"
$ws.Cells.Item(1, 3).Value = "description"

# --- data rows, entered in the same order the workbook author used:
#     name(row2) -> nameForUser(col, rows2-5) -> name(col, rows3-5)
#     -> description(col, rows2-5) -> header nameForUser(B1) ------------
$ws.Cells.Item(2, 1).Value = "Экотуристы"

$ws.Cells.Item(2, 2).Value = "Любитель природы"
$ws.Cells.Item(3, 2).Value = "Охотник за адреналином"
$ws.Cells.Item(4, 2).Value = "Исследователь культур"
$ws.Cells.Item(5, 2).Value = "Ценитель отдыха"

$ws.Cells.Item(3, 1).Value = "Экстремальные туристы"
$ws.Cells.Item(4, 1).Value = "Этнографические туристы"
$ws.Cells.Item(5, 1).Value = "Рекреационные туристы"

$ws.Cells.Item(2, 3).Value = "Путешествуете ради экологии и живописных мест"
$ws.Cells.Item(3, 3).Value = "Ищете драйв, экстрим и новые высоты"
$ws.Cells.Item(4, 3).Value = "Углубляетесь в традиции и быт народов"
$ws.Cells.Item(5, 3).Value = "Наслаждаетесь комфортом, спокойствием и природой"

$ws.Cells.Item(1, 2).Value = "nameForUser"

# --- column widths (best-fit like; inputs chosen so the engine's width
#     quantisation lands on the value closest to the original bestFit
#     widths of 25 / 24.42578125 / 52.28515625) ---------------------------
$ws.Columns.Item(1).ColumnWidth = 24.17
$ws.Columns.Item(2).ColumnWidth = 23.65
$ws.Columns.Item(3).ColumnWidth = 51.5

# --- leave the selection where the author left it ----------------------
$ws.Range("C7").Select() | Out-Null
